$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR input data")
Write-Host $ws.Name
